$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 157; existing rows 157.. shift down to 158..
$ws.Rows("157:157").Insert()

# Populate the new row 157 with its data
$ws.Cells.Item(157, 1).Value = 7
$ws.Cells.Item(157, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(157, 3).Value = "Ñuble"
$ws.Cells.Item(157, 4).Value = 44627
$ws.Cells.Item(157, 5).Value = 16
$ws.Cells.Item(157, 6).Value = 100112043
$ws.Cells.Item(157, 7).Value = "Pepino ensalada"
$ws.Cells.Item(157, 8).Value = "Sin especificar"
$ws.Cells.Item(157, 9).Value = "Primera"
$ws.Cells.Item(157, 10).Value = 100
$ws.Cells.Item(157, 11).Value = 14000
$ws.Cells.Item(157, 12).Value = 15000
$ws.Cells.Item(157, 13).Value = 14500
$ws.Cells.Item(157, 14).Value = "$/caja 80 unidades"
$ws.Cells.Item(157, 15).Value = "Región del Maule"
$ws.Cells.Item(157, 16).Value = 181
$ws.Cells.Item(157, 17).Value = 80
$ws.Cells.Item(157, 18).Value = "Hortaliza"
